$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.003.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.238.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.91%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.621'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.27'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.32%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.47%  '
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.33%  '
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.554.69'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.243.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.814'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.918.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.52%  '
$ws.Range("E20").Value = '  +3.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("E23").Value = '  +4.17%  '
$ws.Range("E24").Value = '  +14.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '230.82'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.90'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  -5.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.37'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +25.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.43%  '
$ws.Range("E32").Value = '  -2.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.32'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.47%  '
$ws.Range("E34").Value = '  +3.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.76%  '
$ws.Range("E36").Value = '  +11.29%  '
$ws.Range("E37").Value = '  +1.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.40'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0328'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +14.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.04%  '
$ws.Range("E41").Value = '  +2.73%  '
$ws.Range("E42").Value = '  +2.10%  '
$ws.Range("E43").Value = '  +6.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.02%  '
$ws.Range("E47").Value = '  +1.28%  '
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("E49").Value = '  +3.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.440'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +15.42%  '
$ws.Range("E51").Value = '  +1.09%  '
